$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.196.01"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "'2.600.49"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'540.96"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "'140.99"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "'3.062.27"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'59.136.90"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'2.621.85"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'341.67"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "'10.10"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'67.60"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "  +8.32%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "'149.39"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.834"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "'0.815"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'273.90"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "'10.73"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'0.0956"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'1.938.76"
$ws.Range("D49").Value = "'18.44"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'111.30"
$ws.Range("E51").Value = "  -1.96%  "
